# Change the sort order of the "Host" column values on the Terminology sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Terminology")

# The sheet is protected, so it must be unprotected before the cells can be
# edited, and re-protected afterwards to restore the original behavior.
$ws.Unprotect()

$ws.Range("A2").Value = "human (Homo sapiens)"
$ws.Range("A3").Value = "mouse (Mus musculus)"
$ws.Range("A4").Value = "chicken (Gallus gallus)"
$ws.Range("A5").Value = "llama (Lama glama)"

$ws.Protect()
